$d = $word.ActiveDocument

# Collapse the split "Fall 2022 INFO-233 Group Project – Phase " / "I" / "I – " / "Coding"
# runs in the title paragraph into a single run reading
# "INFO-233 Group Project – Phase II – Coding" (dropping the "Fall 2022 " prefix).
$d.Content.Find.Execute("Fall 2022 INFO-233", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "INFO-233", 2)
